$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$notes = $s.NotesPage
$shp = $notes.Shapes.AddPlaceholder(2)
$text = "Add separate legend for P2 amp`r`rChange dots for different levels of fix`r`rMove P2 amp up, truncate axis so they don't overlap`r`rChange color"
$shp.TextFrame.TextRange.Text = $text
Write-Host "Paragraphs: $($shp.TextFrame.TextRange.Paragraphs().Count)"
